$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 302.07693
$ws.Range("I12").Value = 302.07693
$ws.Range("K12").Value = 302.07693
$ws.Range("M12").Value = -132.07693
$ws.Range("H19").Value = 523.125
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 523.125
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 523.125
$ws.Range("N19").Value = -873.125
$ws.Range("M19").ClearContents()
$ws.Range("H74").Value = 6381.6665
$ws.Range("I74").Value = 5947.5
$ws.Range("K74").Value = 5947.5
$ws.Range("M74").Value = -5011.5
$ws.Range("H77").Value = 6381.6665
$ws.Range("I77").Value = 5947.5
$ws.Range("K77").Value = 29737.5
$ws.Range("M77").Value = -25057.5
$ws.Range("H100").Value = 95199.31
$ws.Range("I100").Value = 168788.17
$ws.Range("K100").Value = 168788.17
$ws.Range("M100").Value = -168247.17
$ws.Range("H113").Value = 9299.091
$ws.Range("J113").Value = 3799
$ws.Range("L113").Value = 3799
$ws.Range("N113").Value = -10307
$ws.Range("H137").Value = 7345.4683
$ws.Range("I137").Value = 9220.529
$ws.Range("K137").Value = 27661.587
$ws.Range("M137").Value = -25111.587

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3079.8845
$ws.Range("I2").Value = 4036.6428
$ws.Range("J2").Value = 1963.6666
$ws.Range("K2").Value = 4036.6428
$ws.Range("L2").Value = 1963.6666
$ws.Range("M2").Value = -3923.6428
$ws.Range("N2").Value = -2189.6666
$ws.Range("H61").Value = 11030.131
$ws.Range("I61").Value = 13621
$ws.Range("J61").Value = 6999.8887
$ws.Range("K61").Value = 13621
$ws.Range("L61").Value = 6999.8887
$ws.Range("M61").Value = -13409
$ws.Range("N61").Value = -7423.8887
$ws.Range("H116").Value = 3079.8845
$ws.Range("I116").Value = 4036.6428
$ws.Range("J116").Value = 1963.6666
$ws.Range("K116").Value = 4036.6428
$ws.Range("L116").Value = 1963.6666
$ws.Range("M116").Value = -1742.6428
$ws.Range("N116").Value = -6551.6666
$ws.Range("H132").Value = 2706.92
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 5962.25
$ws.Range("K132").Value = 3525
$ws.Range("L132").Value = 17886.75
$ws.Range("M132").Value = -995
$ws.Range("N132").Value = -22946.75
$ws.Range("H136").Value = 11030.131
$ws.Range("I136").Value = 13621
$ws.Range("J136").Value = 6999.8887
$ws.Range("K136").Value = 40863
$ws.Range("L136").Value = 20999.6661
$ws.Range("M136").Value = -38313
$ws.Range("N136").Value = -26099.6661

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3079.8845
$ws.Range("I3").Value = 4036.6428
$ws.Range("J3").Value = 1963.6666
$ws.Range("K3").Value = 4036.6428
$ws.Range("L3").Value = 1963.6666
$ws.Range("M3").Value = -3922.6428
$ws.Range("N3").Value = -2191.6666
$ws.Range("H20").Value = 2514.2903
$ws.Range("I20").Value = 1841.421
$ws.Range("J20").Value = 3579.6667
$ws.Range("K20").Value = 1841.421
$ws.Range("L20").Value = 3579.6667
$ws.Range("M20").Value = -1594.421
$ws.Range("N20").Value = -4073.6667
$ws.Range("H99").Value = 10617.742
$ws.Range("I99").Value = 11358
$ws.Range("K99").Value = 11358
$ws.Range("M99").Value = -9860

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 334149.34
$ws.Range("I16").Value = 999
$ws.Range("J16").Value = 500724.5
$ws.Range("K16").Value = 999
$ws.Range("L16").Value = 500724.5
$ws.Range("M16").Value = -712
$ws.Range("N16").Value = -501298.5
$ws.Range("H31").Value = 8894.206
$ws.Range("I31").Value = 11972.5625
$ws.Range("K31").Value = 11972.5625
$ws.Range("M31").Value = -11677.5625
$ws.Range("H34").Value = 8894.206
$ws.Range("I34").Value = 11972.5625
$ws.Range("K34").Value = 11972.5625
$ws.Range("M34").Value = -11770.5625
$ws.Range("H47").Value = 46999
$ws.Range("J47").Value = 46999
$ws.Range("L47").Value = 46999
$ws.Range("N47").Value = -48131
$ws.Range("H58").Value = 3952.36
$ws.Range("I58").Value = 4367.25
$ws.Range("K58").Value = 4367.25
$ws.Range("M58").Value = -4164.25
$ws.Range("H113").Value = 334149.34
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 500724.5
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 500724.5
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -505064.5
$ws.Range("H132").Value = 1089.6757
$ws.Range("I132").Value = 994
$ws.Range("K132").Value = 2982
$ws.Range("M132").Value = -452
$ws.Range("H136").Value = 3952.36
$ws.Range("I136").Value = 4367.25
$ws.Range("K136").Value = 13101.75
$ws.Range("M136").Value = -10551.75
$ws.Range("H141").Value = 184641.73
$ws.Range("J141").Value = 196924.78
$ws.Range("L141").Value = 196924.78
$ws.Range("N141").Value = -207284.78

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 183.09091
$ws.Range("I33").Value = 93
$ws.Range("J33").Value = 234.57143
$ws.Range("K33").Value = 558
$ws.Range("L33").Value = 1407.42858
$ws.Range("M33").Value = -275
$ws.Range("N33").Value = -1973.42858
$ws.Range("H52").Value = 12671.556
$ws.Range("J52").Value = 12671.556
$ws.Range("L52").Value = 38014.66800000001
$ws.Range("N52").Value = -38546.66800000001
$ws.Range("H131").Value = 2653.3093
$ws.Range("J131").Value = 2014.1548
$ws.Range("L131").Value = 6042.4644
$ws.Range("N131").Value = -16122.4644

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 13406.417
$ws.Range("I70").Value = 12975.5
$ws.Range("J70").Value = 13621.875
$ws.Range("K70").Value = 12975.5
$ws.Range("L70").Value = 13621.875
$ws.Range("M70").Value = -12705.5
$ws.Range("N70").Value = -14161.875
$ws.Range("H73").Value = 13406.417
$ws.Range("I73").Value = 12975.5
$ws.Range("J73").Value = 13621.875
$ws.Range("K73").Value = 12975.5
$ws.Range("L73").Value = 13621.875
$ws.Range("M73").Value = -12039.5
$ws.Range("N73").Value = -15493.875
$ws.Range("H107").Value = 433.57144
$ws.Range("I107").Value = 465
$ws.Range("K107").Value = 465
$ws.Range("M107").Value = 1455
$ws.Range("H113").Value = 7875.1113
$ws.Range("I113").Value = 10413.5
$ws.Range("J113").Value = 2798.3333
$ws.Range("K113").Value = 10413.5
$ws.Range("L113").Value = 2798.3333
$ws.Range("M113").Value = -8243.5
$ws.Range("N113").Value = -7138.3333
$ws.Range("H122").Value = 8349.281999999999
$ws.Range("I122").Value = 6221
$ws.Range("K122").Value = 18663
$ws.Range("M122").Value = -16213

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2029.7931
$ws.Range("J46").Value = 2728.1052
$ws.Range("L46").Value = 2728.1052
$ws.Range("N46").Value = -3104.1052
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H132").Value = 516398.06
$ws.Range("J132").Value = 3660
$ws.Range("L132").Value = 10980
$ws.Range("N132").Value = -16040

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 39149.832
$ws.Range("J75").Value = 45000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46872
$ws.Range("H78").Value = 39149.832
$ws.Range("J78").Value = 45000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -144360
$ws.Range("H132").Value = 7669
$ws.Range("I132").Value = 8991.312
$ws.Range("J132").Value = 4693.8
$ws.Range("K132").Value = 26973.936
$ws.Range("L132").Value = 14081.4
$ws.Range("M132").Value = -24443.936
$ws.Range("N132").Value = -19141.4
$ws.Range("H136").Value = 382142.88
$ws.Range("I136").Value = 485904.94
$ws.Range("J136").Value = 13211.111
$ws.Range("K136").Value = 1457714.82
$ws.Range("L136").Value = 39633.333
$ws.Range("M136").Value = -1455164.82
$ws.Range("N136").Value = -44733.333
$ws.Range("H140").Value = 89166.5
$ws.Range("J140").Value = 89166.5
$ws.Range("L140").Value = 89166.5
$ws.Range("N140").Value = -99526.5
